$wb = $excel.ActiveWorkbook

# --- Sheet "About" (sheet1) ---
$about = $wb.Worksheets.Item("About")

# Update the existing note text (row 10) and add additional note lines (rows 12-15)
$about.Range("A10").Value = "mandate in India for any vehicle type."
$about.Range("A12").Value = "In BAU, there are only policies for incentivising demand "
$about.Range("A13").Value = "creation such as the Faster Adoption and Manufacturing of"
$about.Range("A14").Value = "Hybrid and Electric Vehicles (FAME) scheme which provides"
$about.Range("A15").Value = "subsidies for purchase of different EV types."

# Update selection to reflect the newly added rows
$about.Range("A9:A15").Select()

# --- Sheets "BMRESP-passenger" and "BMRESP-freight" (sheet2 / sheet3) ---
$sheetNames = @("BMRESP-passenger", "BMRESP-freight")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $cell = $ws.Range("A1")
    $cell.Value = "Sales Percentage (dimensionless)"
    $cell.Font.Bold = $true
    $cell.WrapText = $true
    $ws.Rows.Item(1).RowHeight = 30
}
